$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")

# --- "总计" sheet: the existing "2020-Q4" summary row becomes "2022-Q3" ---
# (new quarter's totals), and a new row is appended below it holding the
# original "2020-Q4" totals that used to live there.
$total.Range("A2:D2").Copy($total.Range("A3"))
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2020-Q4"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.01

$total.Range("B2").Value = "2022-Q3"
$total.Range("D2").Value = 0.06

# --- Insert the new "2022-Q3" sheet right before the "2020-Q4" sheet ---
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("2020-Q4"))
$newSheet.Name = "2022-Q3"

# Header row - reuse the bold/bordered header style from "总计"
$total.Range("B1").Copy($newSheet.Range("B1:H1"))
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row - reuse the style from 总计's A2 for the leading index cell
$total.Range("A2").Copy($newSheet.Range("A2"))
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'486002"
$newSheet.Range("C2").Value = "工银全球精选股票（QDII）"
$newSheet.Range("D2").Value = "'3.72"
$newSheet.Range("E2").Value = "'93.69"
$newSheet.Range("F2").Value = "'1.52"
$newSheet.Range("G2").Value = "'0.0565"
$newSheet.Range("H2").Value = 9

# Restore the original tab-selection onto the (now third) "2020-Q4" sheet.
$wb.Worksheets.Item("2020-Q4").Select()
